# GTOR Engine Dyno.docx — documentation edit
#
# 1) Split the run "The load cell" into "The load " + "cell", wrapping
#    "cell" with gramStart/gramEnd proof-error markers (as Word's
#    grammar checker would do), by surgically re-inserting OOXML over
#    just the "cell" substring so the surrounding runs are preserved.
# 2) Remove the two stray empty paragraphs that sit between the
#    "3.0 Software Theory of Operation" heading and the "Arduino"
#    paragraph.

$d = $word.ActiveDocument

# --- 1. "The load cell" -> "The load " / proofErr(gramStart) / "cell" / proofErr(gramEnd) ---

$prefix = "The load "

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("The load cell")) {
        $cellStart = $p.Range.Start + $prefix.Length
        $cellEnd = $cellStart + 4
        $rngCell = $d.Range($cellStart, $cellEnd)

        $cellXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:proofErr w:type="gramStart"/><w:r><w:t>cell</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
        $rngCell.InsertXML($cellXml)
        break
    }
}

# --- 2. Delete the two empty paragraphs after "3.0 Software Theory of Operation" ---

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("3.0 Software Theory of Operation")) {
        $afterHeading = $p.Range.End
        break
    }
}

$d.Range($afterHeading, $afterHeading + 1).Delete()
$d.Range($afterHeading, $afterHeading + 1).Delete()
